# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns, and
# fix the LidoDAOToken / EthereumClassic row ordering (rows 27 & 28 swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Columns D/E store plain display text (e.g. "1.000", "27.714.97") even
    # though many look numeric. Force text storage, then drop back to the
    # default (unstyled) cell style so no stray formatting is left behind.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "27.721.54"
$ws.Range("E2").Value = "  -0.11%  "

Set-TextValue "D3" "1.900.31"
$ws.Range("E3").Value = "  +0.17%  "

Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.21%  "

Set-TextValue "D5" "311.79"
$ws.Range("E5").Value = "  -0.08%  "

Set-TextValue "D6" "1.000"
$ws.Range("E6").Value = "  +0.09%  "

Set-TextValue "D7" "0.5206"
$ws.Range("E7").Value = "  +5.57%  "

Set-TextValue "D8" "0.3782"
$ws.Range("E8").Value = "  -0.31%  "

Set-TextValue "D9" "0.07243"
$ws.Range("E9").Value = "  -1.16%  "

Set-TextValue "D10" "21.20"
$ws.Range("E10").Value = "  +2.74%  "

Set-TextValue "D11" "0.9025"
$ws.Range("E11").Value = "  -1.01%  "

Set-TextValue "D12" "0.07640"
$ws.Range("E12").Value = "  +0.20%  "

Set-TextValue "D13" "1.928.51"
$ws.Range("E13").Value = "  +1.76%  "

Set-TextValue "D14" "5.447"
$ws.Range("E14").Value = "  -0.35%  "

Set-TextValue "D15" "92.14"
$ws.Range("E15").Value = "  +1.06%  "

Set-TextValue "D16" "1.002"
$ws.Range("E16").Value = "  +0.19%  "

Set-TextValue "D17" "0.000008708"
$ws.Range("E17").Value = "  -0.29%  "

Set-TextValue "D18" "1.000"
$ws.Range("E18").Value = "  +0.09%  "

Set-TextValue "D19" "27.773.11"
$ws.Range("E19").Value = "  +0.06%  "

Set-TextValue "D20" "14.46"
$ws.Range("E20").Value = "  -0.16%  "

Set-TextValue "D21" "5.137"
$ws.Range("E21").Value = "  +0.27%  "

Set-TextValue "D22" "2.129.32"
$ws.Range("E22").Value = "  -0.19%  "

Set-TextValue "D23" "10.84"
$ws.Range("E23").Value = "  +0.75%  "

Set-TextValue "D24" "6.617"
$ws.Range("E24").Value = "  -0.63%  "

Set-TextValue "D25" "153.27"
$ws.Range("E25").Value = "  -0.49%  "

Set-TextValue "D26" "1.867"
$ws.Range("E26").Value = "  +0.92%  "

# Rows 27 & 28 swap places: EthereumClassic/LidoDAOToken ordering flips,
# each keeping its own refreshed price + volume figures.
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D27" "2.161"
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "18.28"
$ws.Range("E28").Value = "  -0.65%  "

Set-TextValue "D29" "114.39"
$ws.Range("E29").Value = "  -0.96%  "

Set-TextValue "D30" "4.840"
$ws.Range("E30").Value = "  -0.89%  "

Set-TextValue "D31" "0.09076"
$ws.Range("E31").Value = "  +1.57%  "

Set-TextValue "D32" "3.188"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("E33").Value = "  +4.37%  "

Set-TextValue "D34" "1.227"
$ws.Range("E34").Value = "  -0.20%  "

Set-TextValue "D35" "0.7796"
$ws.Range("E35").Value = "  +1.53%  "

Set-TextValue "D36" "0.02089"
$ws.Range("E36").Value = "  +2.23%  "

Set-TextValue "D37" "2.589"
$ws.Range("E37").Value = "  +0.83%  "

Set-TextValue "D38" "3.072"
$ws.Range("E38").Value = "  +2.74%  "

Set-TextValue "D39" "1.093"
$ws.Range("E39").Value = "  -0.66%  "

Set-TextValue "D40" "0.5553"
$ws.Range("E40").Value = "  +0.94%  "

Set-TextValue "D41" "0.05294"
$ws.Range("E41").Value = "  -0.04%  "

Set-TextValue "D42" "6.723"
$ws.Range("E42").Value = "  -2.68%  "

Set-TextValue "D43" "116.05"
$ws.Range("E43").Value = "  +3.35%  "

Set-TextValue "D44" "8.507"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("E45").Value = "  -0.36%  "

Set-TextValue "D46" "0.4815"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("E47").Value = "  -1.42%  "

Set-TextValue "D48" "0.9995"
$ws.Range("E48").Value = "  +0.04%  "

Set-TextValue "D49" "1.613"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("E50").Value = "  -1.14%  "

$ws.Range("E51").Value = "  -0.79%  "
